# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Column D ("Price") cells are stored as literal text (e.g. "40.724.29" uses
# dots as thousands separators), and column E ("Volume(1h)") cells are padded
# percentage strings. A leading apostrophe is used where needed so Excel keeps
# plain-decimal-looking price text (e.g. "311.07") as text instead of a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.724.29"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "2.378.13"
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'311.07"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "'87.64"
$ws.Range("E6").Value = "  -6.04%  "
$ws.Range("E7").Value = "  -4.41%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -4.43%  "
$ws.Range("D10").Value = "'0.0845"
$ws.Range("E10").Value = "  -4.58%  "
$ws.Range("D11").Value = "'30.81"
$ws.Range("E11").Value = "  -7.48%  "
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").Value = "2.745.39"
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("D14").Value = "'6.57"
$ws.Range("E14").Value = "  -5.20%  "
$ws.Range("D15").Value = "'15.01"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "2.391.82"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "'0.763"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "40.647.58"
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = "  -4.94%  "
$ws.Range("D21").Value = "'68.91"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").Value = "'232.98"
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "  -4.78%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'1.81"
$ws.Range("E26").Value = "  -7.22%  "
$ws.Range("D27").Value = "'23.87"
$ws.Range("E27").Value = "  -5.56%  "
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").Value = "'9.38"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("D30").Value = "'34.06"
$ws.Range("E30").Value = "  -7.80%  "
$ws.Range("D31").Value = "'152.91"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.24"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").Value = "'0.0733"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("E36").Value = "  -2.43%  "
$ws.Range("D37").Value = "'16.01"
$ws.Range("E37").Value = "  -8.15%  "
$ws.Range("D38").Value = "'2.78"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").Value = "'0.0999"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("D40").Value = "'1.72"
$ws.Range("E40").Value = "  -8.32%  "
$ws.Range("D41").Value = "'3.89"
$ws.Range("E41").Value = "  -3.79%  "
$ws.Range("D42").Value = "'2.41"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("D43").Value = "1.961.91"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").Value = "'0.0271"
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("D45").Value = "'17.62"
$ws.Range("E45").Value = "  -9.10%  "
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -8.40%  "
$ws.Range("D48").Value = "2.602.67"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("D50").Value = "'72.56"
$ws.Range("E50").Value = "  -5.76%  "
$ws.Range("D51").Value = "'50.80"
$ws.Range("E51").Value = "  -3.10%  "
